# Add a new "ShowBorder" bool column to the Tile table (table1 / "表1"),
# mirroring the existing Id/Name/Cname/... columns: header row (row1),
# type/localized-header rows (rows 2-3), then "true" for every data row
# (rows 4-12). Column H, table grows from A1:G13 to A1:H13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Growing the table via ListColumns.Add() extends tableColumns, the table
# ref, the autoFilter ref and (once the header cell gets a value below)
# the sheet dimension / row spans too - exactly like typing into the
# first blank column to the right of an Excel table.
$col = $lo.ListColumns.Add()

# Row 1: column header
$col.Range.Cells.Item(1, 1).Value = "ShowBorder"

# Row 2: field type
$ws.Range("H2").Value = "bool"

# Row 3: Chinese field label
$ws.Range("H3").Value = "显示边框"

# Rows 4-12: data values - every tile row defaults to ShowBorder = true.
# A leading apostrophe forces these to be stored as text "true" (matching
# the existing E/F-column text cells) instead of being auto-coerced to a
# native Excel boolean.
for ($r = 4; $r -le 12; $r++) {
    $ws.Cells.Item($r, 8).Value = "'true"
}

# Match formatting of the sibling columns: row2/row3 header styling comes
# from column G, the data rows' text-number-format styling comes from
# column F (same style used for the existing CanMove/Icon text columns).
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)

$ws.Range("F4:F12").Copy()
$ws.Range("H4:H12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Selection ends up on the newly-filled data column, same as after typing
# the new values in by hand.
[void]$ws.Range("H4:H12").Select()
